$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update LAST SCRAPE DATE column (F) for all data rows: 2019-03-07 -> 2019-03-12
# (leading apostrophe forces text so Excel doesn't reinterpret the string as a
# date serial number; ClearFormats() strips the resulting quote-prefix style
# so the cell stays on the default style, matching every other text cell)
$ws.Range("F2:F70").Value = "'2019-03-12"
$ws.Range("F2:F70").ClearFormats()

# Update GAME NAME (C), GAME NUMBER (D), TOP PRIZES REMAINING (E) for rows whose data moved
# Row 6
$ws.Range("C6").Value = "BETTY BOOP™"
$ws.Range("D6").Value = 471
$ws.Range("E6").Value = 1
# Row 7
$ws.Range("E7").Value = 1972
# Row 8
$ws.Range("C8").Value = "Super 7's"
$ws.Range("D8").Value = 485
$ws.Range("E8").Value = 3
# Row 9
$ws.Range("C9").Value = "10X Payout"
$ws.Range("D9").Value = 490
$ws.Range("E9").Value = 4
# Row 10
$ws.Range("C10").Value = "Cash $100's"
$ws.Range("D10").Value = 425
$ws.Range("E10").Value = 813
# Row 11
$ws.Range("E11").Value = 767
# Row 21
$ws.Range("C21").Value = "Jumbo Bucks"
$ws.Range("D21").Value = 486
$ws.Range("E21").Value = 3
# Row 22
$ws.Range("C22").Value = "20X The Money"
$ws.Range("D22").Value = 435
$ws.Range("E22").Value = 1
# Row 26
$ws.Range("E26").Value = 1251
# Row 28
$ws.Range("E28").Value = 936
# Row 42
$ws.Range("C42").Value = "50X Payout"
$ws.Range("D42").Value = 492
$ws.Range("E42").Value = 4
# Row 43
$ws.Range("C43").Value = "Cash Multiplier"
$ws.Range("D43").Value = 478
$ws.Range("E43").Value = 2
# Row 44
$ws.Range("C44").Value = "Stacks of Cash"
$ws.Range("D44").Value = 465
$ws.Range("E44").Value = 1
# Row 45
$ws.Range("C45").Value = "$100,000 Triple Win"
$ws.Range("D45").Value = 474
$ws.Range("E45").Value = 1
# Row 46
$ws.Range("C46").Value = "Hit $500!"
$ws.Range("D46").Value = 483
$ws.Range("E46").Value = 1781
# Row 52
$ws.Range("C52").Value = "$250,000 Riches"
$ws.Range("D52").Value = 466
# Row 53
$ws.Range("C53").Value = "$200,000 Jackpot"
$ws.Range("D53").Value = 447
# Row 55
$ws.Range("C55").Value = "$200,000 Bonus"
$ws.Range("D55").Value = 489
$ws.Range("E55").Value = 3
# Row 56
$ws.Range("C56").Value = "Big Multiplier"
$ws.Range("D56").Value = 470
$ws.Range("E56").Value = 2
# Row 58
$ws.Range("C58").Value = "100X Payout"
$ws.Range("D58").Value = 493
$ws.Range("E58").Value = 4
# Row 59
$ws.Range("C59").Value = "Cash $1000's"
$ws.Range("D59").Value = 428
$ws.Range("E59").Value = 61
# Row 60
$ws.Range("E60").Value = 168
# Row 62
$ws.Range("E62").Value = 819
# Row 69
$ws.Range("C69").Value = "$50K Blowout"
$ws.Range("D69").Value = 461
$ws.Range("E69").Value = 29
# Row 70
$ws.Range("C70").Value = "200X Payout"
$ws.Range("D70").Value = 494
$ws.Range("E70").Value = 4
